$d = $word.ActiveDocument

# Locate the paragraph that currently starts with "Perseusz " (the one whose
# run-level content is being collapsed into a single clean run that also
# absorbs/duplicates the opening sentence naming the constellation).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Perseusz")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph (starts with 'Perseusz')"
}

$newText = "Uczestniczysz w ogólnoświatowym przedsięwzięciu, którego celem jest obserwacja i odnotowanie najsłabszych widocznych gwiazd w celu zmierzenia zanieczyszczenia światłem w danym miejscu. Poprzez zlokalizowanie i obserwację  Gwiazdozbiór Bliźniąt na nocnym niebie oraz porównanie go do map nieba ludzie z całego świata będą mogli dowiedzieć się jaki wkład światło emitowane przez ich społeczność wnosi do  zanieczyszczenia światłem. To co dodasz do internetowej bazy danych pomoże udokumentować widoczne nocne niebo."

# Capture the paragraph's own attributes (paraId/rsidR/etc.) and its
# existing <w:pPr> so that InsertXML - which replaces the *entire* <w:p> it
# lands on - recreates them unchanged, while every old run and every stray
# <w:proofErr/> marker inside the paragraph gets swept away and replaced by
# one plain run holding the new merged text.
$pPrXml = '<w:pPr><w:pStyle w:val="BasicParagraph"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="-72"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Optima-Regular"/></w:rPr></w:pPr>'
$pAttrs = 'w14:paraId="38C52869" w14:textId="35A5E625" w:rsidR="00292489" w:rsidRDefault="00310DC5" w:rsidP="00C92045"'

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + $pAttrs + '>' + $pPrXml + '<w:r><w:t>' + $newText + '</w:t></w:r></w:p>'

$target.Range.InsertXML($xml)
Write-Output "Paragraph replaced."
